$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H76").Value = 3000.2727
$ws.Range("I76").Value = 3000.2727
$ws.Range("K76").Value = 3000.2727
$ws.Range("M76").Value = -2685.2727
$ws.Range("H79").Value = 3000.2727
$ws.Range("I79").Value = 3000.2727
$ws.Range("K79").Value = 3000.2727
$ws.Range("M79").Value = -1908.2727
$ws.Range("H141").Value = 1660.6
$ws.Range("I141").Value = 1071.2941
$ws.Range("K141").Value = 3213.8823
$ws.Range("M141").Value = 1966.1177

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 3062.3809
$ws.Range("I63").Value = 2513.75
$ws.Range("J63").Value = 3400
$ws.Range("K63").Value = 2513.75
$ws.Range("L63").Value = 3400
$ws.Range("M63").Value = -1827.75
$ws.Range("N63").Value = -4772
$ws.Range("H66").Value = 3062.3809
$ws.Range("I66").Value = 2513.75
$ws.Range("J66").Value = 3400
$ws.Range("K66").Value = 12568.75
$ws.Range("L66").Value = 17000
$ws.Range("M66").Value = -9136.75
$ws.Range("N66").Value = -23864
$ws.Range("H74").Value = 22695.83
$ws.Range("I74").Value = 28095.838
$ws.Range("J74").Value = 2715.8
$ws.Range("K74").Value = 28095.838
$ws.Range("L74").Value = 2715.8
$ws.Range("M74").Value = -27221.838
$ws.Range("N74").Value = -4463.8
$ws.Range("H77").Value = 22695.83
$ws.Range("I77").Value = 28095.838
$ws.Range("J77").Value = 2715.8
$ws.Range("K77").Value = 140479.19
$ws.Range("L77").Value = 13579
$ws.Range("M77").Value = -136111.19
$ws.Range("N77").Value = -22315
$ws.Range("H122").Value = 1956.75
$ws.Range("I122").Value = 1595.8
$ws.Range("K122").Value = 4787.4
$ws.Range("M122").Value = -2337.4
$ws.Range("H123").Value = 21000
$ws.Range("J123").Value = 21000
$ws.Range("L123").Value = 21000
$ws.Range("N123").Value = -30800
$ws.Range("H124").Value = 23221.416
$ws.Range("J124").Value = 23221.416
$ws.Range("L124").Value = 23221.416
$ws.Range("N124").Value = -33041.416
$ws.Range("H127").Value = 33779
$ws.Range("J127").Value = 33779
$ws.Range("L127").Value = 33779
$ws.Range("N127").Value = -43699
$ws.Range("H128").Value = 34113.8
$ws.Range("J128").Value = 34113.8
$ws.Range("L128").Value = 34113.8
$ws.Range("N128").Value = -44073.8
$ws.Range("H129").Value = 35000
$ws.Range("J129").Value = 35000
$ws.Range("L129").Value = 35000
$ws.Range("N129").Value = -45000
$ws.Range("H130").Value = 29775.555
$ws.Range("J130").Value = 29775.555
$ws.Range("L130").Value = 29775.555
$ws.Range("N130").Value = -39815.555
$ws.Range("H131").Value = 34318.184
$ws.Range("J131").Value = 34318.184
$ws.Range("L131").Value = 34318.184
$ws.Range("N131").Value = -44398.184

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1694.6428
$ws.Range("I86").Value = 1516.5
$ws.Range("J86").Value = 2347.8333
$ws.Range("K86").Value = 1516.5
$ws.Range("L86").Value = 2347.8333
$ws.Range("M86").Value = -393.5
$ws.Range("N86").Value = -4593.8333
$ws.Range("H89").Value = 1694.6428
$ws.Range("I89").Value = 1516.5
$ws.Range("J89").Value = 2347.8333
$ws.Range("K89").Value = 7582.5
$ws.Range("L89").Value = 11739.1665
$ws.Range("M89").Value = -1966.5
$ws.Range("N89").Value = -22971.1665
$ws.Range("H123").Value = 25000
$ws.Range("J123").Value = 25000
$ws.Range("L123").Value = 25000
$ws.Range("N123").Value = -34800
$ws.Range("H126").Value = 35000
$ws.Range("J126").Value = 35000
$ws.Range("L126").Value = 35000
$ws.Range("N126").Value = -44880
$ws.Range("H127").Value = 0
$ws.Range("J127").Value = 0
$ws.Range("L127").Value = 0
$ws.Range("N127").ClearContents()
$ws.Range("H129").Value = 35559.082
$ws.Range("I129").Value = 30709
$ws.Range("J129").Value = 36000
$ws.Range("K129").Value = 30709
$ws.Range("L129").Value = 36000
$ws.Range("M129").Value = -25709
$ws.Range("N129").Value = -46000
$ws.Range("H131").Value = 30000
$ws.Range("J131").Value = 30000
$ws.Range("L131").Value = 30000
$ws.Range("N131").Value = -40080

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 7938341
$ws.Range("I31").Value = 1184.7941
$ws.Range("J31").Value = 17243972
$ws.Range("K31").Value = 1184.7941
$ws.Range("L31").Value = 17243972
$ws.Range("M31").Value = -889.7941000000001
$ws.Range("N31").Value = -17244562
$ws.Range("H34").Value = 7938341
$ws.Range("I34").Value = 1184.7941
$ws.Range("J34").Value = 17243972
$ws.Range("K34").Value = 1184.7941
$ws.Range("L34").Value = 17243972
$ws.Range("M34").Value = -982.7941000000001
$ws.Range("N34").Value = -17244376
$ws.Range("H132").Value = 1374232.5
$ws.Range("I132").Value = 2668.6667
$ws.Range("J132").Value = 6174706
$ws.Range("K132").Value = 8006.000100000001
$ws.Range("L132").Value = 18524118
$ws.Range("M132").Value = -5476.000100000001
$ws.Range("N132").Value = -18529178

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H100").Value = 3600
$ws.Range("J100").Value = 3600
$ws.Range("L100").Value = 10800
$ws.Range("N100").Value = -12422
$ws.Range("H109").Value = 2028.4
$ws.Range("I109").Value = 1051.909
$ws.Range("J109").Value = 3221.889
$ws.Range("K109").Value = 3155.727
$ws.Range("L109").Value = 9665.667000000001
$ws.Range("M109").Value = -2115.727
$ws.Range("N109").Value = -11745.667
$ws.Range("H115").Value = 2024.8889
$ws.Range("I115").Value = 996
$ws.Range("J115").Value = 2848
$ws.Range("K115").Value = 2988
$ws.Range("L115").Value = 8544
$ws.Range("M115").Value = -1813
$ws.Range("N115").Value = -10894

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 257638.75
$ws.Range("I122").Value = 338333.34
$ws.Range("J122").Value = 15555
$ws.Range("K122").Value = 1015000.02
$ws.Range("L122").Value = 46665
$ws.Range("M122").Value = -1012550.02
$ws.Range("N122").Value = -51565
$ws.Range("H123").Value = 12628.866
$ws.Range("J123").Value = 12628.866
$ws.Range("L123").Value = 12628.866
$ws.Range("N123").Value = -17528.866
$ws.Range("H124").Value = 25600
$ws.Range("J124").Value = 25600
$ws.Range("L124").Value = 25600
$ws.Range("N124").Value = -35420
$ws.Range("H125").Value = 28494.5
$ws.Range("J125").Value = 28494.5
$ws.Range("L125").Value = 28494.5
$ws.Range("N125").Value = -33414.5
$ws.Range("H126").Value = 2551
$ws.Range("I126").Value = 3058.2
$ws.Range("J126").Value = 2128.3333
$ws.Range("K126").Value = 9174.599999999999
$ws.Range("L126").Value = 6384.999899999999
$ws.Range("M126").Value = -6704.599999999999
$ws.Range("N126").Value = -11324.9999
$ws.Range("H127").Value = 0
$ws.Range("J127").Value = 0
$ws.Range("L127").Value = 0
$ws.Range("N127").ClearContents()
$ws.Range("H130").Value = 118000
$ws.Range("J130").Value = 118000
$ws.Range("L130").Value = 118000
$ws.Range("N130").Value = -128040
$ws.Range("H131").Value = 28325
$ws.Range("J131").Value = 28325
$ws.Range("L131").Value = 28325
$ws.Range("N131").Value = -38405
$ws.Range("H132").Value = 2566763.2
$ws.Range("I132").Value = 2595.2222
$ws.Range("K132").Value = 7785.6666
$ws.Range("M132").Value = -5255.6666

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H127").Value = 34155.453
$ws.Range("J127").Value = 34155.453
$ws.Range("L127").Value = 34155.453
$ws.Range("N127").Value = -44075.453

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H125").Value = 21907.188
$ws.Range("J125").Value = 21907.188
$ws.Range("L125").Value = 21907.188
$ws.Range("N125").Value = -31747.188
$ws.Range("H132").Value = 3436.394
$ws.Range("I132").Value = 4332.4736
$ws.Range("J132").Value = 2220.2856
$ws.Range("K132").Value = 12997.4208
$ws.Range("L132").Value = 6660.8568
$ws.Range("M132").Value = -10467.4208
$ws.Range("N132").Value = -11720.8568
